$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cells (row 1): dates for the three newly appended quarters ---
# Copy the formatting of the last existing header cell (BF1) onto the new ones
# so BG1:BI1 keep the bold / centered / bordered header style.
$ws.Range("BF1").Copy()
$ws.Range("BG1:BI1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("BG1").Value = "31/12/2023"
$ws.Range("BH1").Value = "31/03/2024"
$ws.Range("BI1").Value = "30/06/2024"

# --- 2. New data columns BG, BH, BI for rows 2-80 ---
# Each tuple is (row, BG-value, BH-value, BI-value)
$data = @(
    @(2, 3076137.984, 2936739.072, 2962825.984),
    @(3, 138882, 107087, 161231.008),
    @(4, 17600, 14321, 78182),
    @(5, 5180, 4781, 4868),
    @(6, 32930, 26745, 30600),
    @(7, 0, 0, 0),
    @(8, 0, 0, 0),
    @(9, 20300, 17452, 12961),
    @(10, 1620, 1024, 1028),
    @(11, 61252, 42764, 33592),
    @(12, 315255.008, 226976, 217099.008),
    @(13, 0, 0, 0),
    @(14, 0, 0, 0),
    @(15, 312, 315, 21),
    @(16, 0, 0, 0),
    @(17, 0, 0, 0),
    @(18, 0, 0, 0),
    @(19, 4799, 5102, 3529),
    @(20, 0, 0, 0),
    @(21, 0, 0, 0),
    @(22, 0, 0, 0),
    @(23, 2622000.896, 2602675.968, 2584496.128),
    @(24, 0, 0, 0),
    @(25, 0, 0, 0),
    @(26, 3076137.984, 2936739.072, 2962825.984),
    @(27, 437913.984, 473134.016, 506687.008),
    @(28, 7220, 8766, 6101),
    @(29, 63750, 46843, 45795),
    @(30, 15033, 9511, 6807),
    @(31, 172754, 225503.008, 253178),
    @(32, 0, 0, 0),
    @(33, 0, 0, 0),
    @(34, 179156.992, 182511.008, 194806),
    @(35, 0, 0, 0),
    @(36, 0, 0, 0),
    @(37, 1705203.968, 1592500.992, 1571192.064),
    @(38, 901872, 876820.992, 876849.9840000001),
    @(39, 0, 0, 0),
    @(40, 417232, 407500, 428267.008),
    @(41, 0, 0, 0),
    @(42, 0, 0, 0),
    @(43, 386100, 308180, 266075.008),
    @(44, 0, 0, 0),
    @(45, 0, 0, 0),
    @(46, 0, 0, 0),
    @(47, 933020.032, 871104, 884947.008),
    @(48, 4128636.928, 4128636.928, 4128636.928),
    @(49, 1, 1, 1),
    @(50, 0, 0, 0),
    @(51, 0, 0, 0),
    @(52, -3195618.048, -3257533.952, -3243691.008),
    @(53, 0, 0, 0),
    @(54, 0, 0, 0),
    @(55, 0, 0, 0),
    @(56, 0, 0, 0),
    @(59, 50584, 42587, 56022),
    @(60, 104706, -57538, -59112),
    @(61, 155289.984, -14951, -3090),
    @(62, 0, 0, 0),
    @(63, -17934, -14421, -11615),
    @(64, 0, 0, 0),
    @(65, 0, 0, 0),
    @(66, 71126, -1058, 62834),
    @(67, 6878, 0, 0),
    @(68, -31035, -33814, -35947),
    @(69, 757, 512, 1378),
    @(70, -31792.008, -34326, -37325),
    @(74, 184324.992, -64244, 12182),
    @(75, -2193, -1791, -16551),
    @(76, -56276, 4119, 18212),
    @(79, 0, 0, 0),
    @(80, 125856, -61916, 13843)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 59).Value = $entry[1]
    $ws.Cells.Item($r, 60).Value = $entry[2]
    $ws.Cells.Item($r, 61).Value = $entry[3]
}

# --- 3. Rows that are blank "header/section" rows in the original sheet (57, 58,
#        71, 72, 73, 77, 78) gain matching blank cells in BG/BH/BI too. ---
$blankRows = @(57, 58, 71, 72, 73, 77, 78)
foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 59).Value = ""
    $ws.Cells.Item($r, 60).Value = ""
    $ws.Cells.Item($r, 61).Value = ""
}

Write-Output "done"
